# Legacy GSC export data update: the exported "Chart" sheet is a rolling
# 87-day window of daily video-indexing stats. Each refresh drops the
# oldest date (top data row) and everything shifts up by one row, with a
# new date implicitly appended at the bottom (the trailing placeholder
# rows already carry zeroed metrics, so no new row needs to be appended
# here - only the oldest day falls out of the window).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

# Row 1 is the header (Date / No video indexed / Video indexed / Impressions).
# Row 2 is the oldest tracked date (2025-11-18) - drop it so every
# subsequent row shifts up by one, the dimension shrinks by a row, and all
# the shared-string-backed date labels advance by a day.
$ws.Rows.Item(2).Delete()
